# Explore-courses assessment workbook: add a second "grouping" sheet and
# correct a handful of hour totals on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 corrections -----------------------------------------------
$ws1.Range("G4").Value = 100
$ws1.Range("H9").Value = 15
$ws1.Range("H16").Value = 15

# Move the active selection on Sheet1 to G5.
$ws1.Range("G5").Select() | Out-Null

# --- Add Sheet2 (placed right after Sheet1) ----------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Row 9
$ws2.Range("A9").Value = "courses taught by julie zelenski"
$ws2.Range("G9").Value = 0
$ws2.Range("H9").Value = 100

# Row 10 (wrapped, taller row)
$ws2.Range("A10").Value = "courses taught by mehran sahami"
$ws2.Range("A10").WrapText = $true
$ws2.Range("G10").Value = 0
$ws2.Range("H10").Value = 200
$ws2.Rows.Item(10).RowHeight = 48

# Row 11
$ws2.Range("A11").Value = "leon simon"
$ws2.Range("G11").Value = 0
$ws2.Range("H11").Value = 100

# Row 17 (wrapped, taller row)
$ws2.Range("A17").Value = "amelang107b and csre14n"
$ws2.Range("A17").WrapText = $true
$ws2.Range("G17").Value = 0
$ws2.Range("H17").Value = 0
$ws2.Rows.Item(17).RowHeight = 36

# Row 19
$ws2.Range("A19").Value = "math52h and cs105"
$ws2.Range("G19").Value = 0
$ws2.Range("H19").Value = 100

# Row 21
$ws2.Range("A21").Value = "amelang129a and cs109"
$ws2.Range("G21").Value = 0
$ws2.Range("H21").Value = 100

# Rows 25-34 (a second group)
$ws2.Range("A25").Value = "introduction to computing principles"
$ws2.Range("G25").Value = 100
$ws2.Range("H25").Value = 0

$ws2.Range("A26").Value = "the mathematics of the rubik's cube"
$ws2.Range("G26").Value = 0
$ws2.Range("H26").Value = 0

$ws2.Range("A27").Value = "third-year persian, second quarter"
$ws2.Range("G27").Value = 0
$ws2.Range("H27").Value = 0

$ws2.Range("A28").Value = "identity and popular music (femgen 140g, music 140g) , csre140g"
$ws2.Range("G28").Value = 0
$ws2.Range("H28").Value = 100

$ws2.Range("A29").Value = "what is hemispheric"
$ws2.Range("G29").Value = 100
$ws2.Range("H29").Value = 0

$ws2.Range("A30").Value = "first-year hausa"
$ws2.Range("G30").Value = 0
$ws2.Range("H30").Value = 0

$ws2.Range("A31").Value = "first-year hebrew, first quarter (jewishst 101a)"
$ws2.Range("G31").Value = 0
$ws2.Range("H31").Value = 15

$ws2.Range("A32").Value = "comparative fictions of ethnicity (amstud 51q, complit 51q) , csre51q"
$ws2.Range("G32").Value = 0
$ws2.Range("H32").Value = 100

$ws2.Range("A33").Value = "growing up bilingual (chilatst 14n, educ 114n)"
$ws2.Range("G33").Value = 0
$ws2.Range("H33").Value = 0

$ws2.Range("A34").Value = "digital dilemmas"
$ws2.Range("G34").Value = 100
$ws2.Range("H34").Value = 100

$ws2.Range("G35").Formula = "=SUM(G25:G34)"
$ws2.Range("H35").Formula = "=SUM(H25:H34)"

# Rows 70-74 (a third group)
$ws2.Range("A70").Value = "amelang110a"
$ws2.Range("G70").Value = 100
$ws2.Range("H70").Value = 100

$ws2.Range("A71").Value = "math120"
$ws2.Range("G71").Value = 100
$ws2.Range("H71").Value = 100

$ws2.Range("A72").Value = "amelang 144b"
$ws2.Range("G72").Value = 100
$ws2.Range("H72").Value = 60

$ws2.Range("A73").Value = "math 53"
$ws2.Range("G73").Value = 0
$ws2.Range("H73").Value = 15

$ws2.Range("A74").Value = "amelang128b"
$ws2.Range("G74").Value = 100
$ws2.Range("H74").Value = 100

$ws2.Range("G75").Formula = "=SUM(G70:G74)"
$ws2.Range("H75").Formula = "=SUM(H70:H74)"

# Sheet2's own view: scrolled down with F60 selected.
$ws2.Activate() | Out-Null
$ws2.Range("F60").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 44

# Re-activate Sheet1 to match the original "tabSelected" sheet.
$ws1.Activate() | Out-Null
